# Adds ArtisanCommands moveBackground and pidLookahead to the "Commands" sheet
# of the eventbuttons workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# --- Insert "pidLookahead(<int>)" row right after the "pidSource(<int>)" row ---
$ws.Rows.Item(83).Insert()
$ws.Range("B83").Value = "pidLookahead(<int>)"
$ws.Range("C83").Value = "sets the PID lookahead"
$ws.Rows.Item(83).RowHeight = 13.8

# --- Insert "moveBackground(<direction>,<int>)" row right after the "alarmset(<as>)" row ---
$ws.Rows.Item(95).Insert()
$ws.Range("B95").Value = "moveBackground(<direction>,<int>)"
$ws.Range("C95").Value = "moves the background profile the indicated number of steps towards <direction>, with <direction> one of up, down, left, right"
$ws.Rows.Item(95).RowHeight = 13.8

$ws.Range("C95").Select()
